$wb = $excel.ActiveWorkbook

$newGuid = "0cca6f66-9142-4246-ae98-83296c2f2571"
$newHash = "dd50b3fd6cad1ef647c6b9586bdf13f9ca91604f"

# ---------------------------------------------------------------------------
# Sheet "Overview": update file name / path cells and the B2 hyperlink text,
# plus the "Latest HO Xliff Generate Date" timestamp.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-29 11:01:21"

$ovHyperlink = $null
foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$2') { $ovHyperlink = $h }
}
if ($ovHyperlink -ne $null) {
    $ovHyperlink.TextToDisplay = "e2e\$newGuid.md"
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn": update file references, handoff/handback timestamps, clear
# the now-unused "Latest Target File" / "Latest Handback File" columns and
# drop the I2 hyperlink, and resize columns I/J.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-29 11:01:16"
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq '$I$2') { $h.Delete() }
}

$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("J2").Style = "Normal"

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

$zhHyperlink = $null
foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') { $zhHyperlink = $h }
}
if ($zhHyperlink -ne $null) {
    $zhHyperlink.TextToDisplay = "$newGuid.md"
}

# ---------------------------------------------------------------------------
# Sheet "de-de": same treatment.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-29 11:01:21"
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$I$2') { $h.Delete() }
}

$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("J2").Style = "Normal"

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426

$deHyperlink = $null
foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') { $deHyperlink = $h }
}
if ($deHyperlink -ne $null) {
    $deHyperlink.TextToDisplay = "$newGuid.md"
}
